$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append two new test rows (94 and 95), following the same pattern as existing rows
$ws.Range("A94").Value = "Average Length"
$ws.Range("B94").Value = "Test average length"
$ws.Range("C94").Value = "Average_Length_test"

$ws.Range("A95").Value = "Average Recovery"
$ws.Range("B95").Value = "Test average recovery"
$ws.Range("C95").Value = "Average_Recovery_test"

# Update selection to match diff (A100 instead of B100)
$ws.Range("A100").Select()
